$d = $word.ActiveDocument

# 1. "Data: 00/00/2023" -> "Data: 00/00/0000"
$d.Content.Find.Execute("00/00/2023", $false, $false, $false, $false, $false,
                         $true, 1, $false, "00/00/0000", 2)

# 2. "Horário: 00:00 - 00:00 " -> "Horário: 00:00 - 00:00" (drop trailing space)
$d.Content.Find.Execute("00:00 - 00:00 ", $false, $false, $false, $false, $false,
                         $true, 1, $false, "00:00 - 00:00", 2)

# 3. Remove the extra empty paragraph sitting between the "Participantes presentes:"
#    bullet list item and the "Participantes ausentes:" paragraph.
$paras = $d.Paragraphs
for ($i = $paras.Count; $i -ge 1; $i--) {
    $p = $paras.Item($i)
    if ($p.Range.Text -eq "`r" -and $i -lt $paras.Count -and $paras.Item($i + 1).Range.Text -like "Participantes ausentes:*") {
        $p.Range.Delete()
        break
    }
}

# 4. Remove bold from the paragraph mark of the empty paragraph right before
#    "Assuntos discutidos e principais decisões:"
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -eq "`r" -and $i -lt $paras.Count -and $paras.Item($i + 1).Range.Text -like "Assuntos discutidos*") {
        $p.Range.Font.Bold = 0
        break
    }
}
